# Scheduled runner refresh: updated Universalis market-price snapshots and
# recomputed Leve profit figures across the per-crafter-job Sheets.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 59697.332
$ws.Range("I70").Value = 1629.8
$ws.Range("K70").Value = 4889.4
$ws.Range("M70").Value = -4619.4
$ws.Range("H73").Value = 59697.332
$ws.Range("I73").Value = 1629.8
$ws.Range("K73").Value = 4889.4
$ws.Range("M73").Value = -3953.4
$ws.Range("H106").Value = 49165
$ws.Range("I106").Value = 49165
$ws.Range("K106").Value = 49165
$ws.Range("M106").Value = -48534
$ws.Range("H116").Value = 6267.375
$ws.Range("I116").Value = 3787.5
$ws.Range("J116").Value = 8747.25
$ws.Range("K116").Value = 3787.5
$ws.Range("L116").Value = 8747.25
$ws.Range("M116").Value = -345.5
$ws.Range("N116").Value = -15631.25
$ws.Range("H129").Value = 5954.1113
$ws.Range("J129").Value = 6086
$ws.Range("L129").Value = 18258
$ws.Range("N129").Value = -28258
$ws.Range("H132").Value = 2076.0977
$ws.Range("I132").Value = 2006.7778
$ws.Range("K132").Value = 6020.3334
$ws.Range("M132").Value = -3490.3334
$ws.Range("H137").Value = 10418.9375
$ws.Range("J137").Value = 16255.667
$ws.Range("L137").Value = 48767.001
$ws.Range("N137").Value = -53867.001
$ws.Range("H138").Value = 5921.4023
$ws.Range("I138").Value = 6476.5415
$ws.Range("J138").Value = 5691.6895
$ws.Range("K138").Value = 19429.6245
$ws.Range("L138").Value = 17075.0685
$ws.Range("M138").Value = -14289.6245
$ws.Range("N138").Value = -27355.0685

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22849.404
$ws.Range("I32").Value = 12018.857
$ws.Range("K32").Value = 12018.857
$ws.Range("M32").Value = -11731.857
$ws.Range("H61").Value = 2719.5715
$ws.Range("I61").Value = 2719.5715
$ws.Range("K61").Value = 2719.5715
$ws.Range("M61").Value = -2507.5715
$ws.Range("H102").Value = 1833.0416
$ws.Range("I102").Value = 1516.7273
$ws.Range("K102").Value = 1516.7273
$ws.Range("M102").Value = 105.2727
$ws.Range("H136").Value = 2719.5715
$ws.Range("I136").Value = 2719.5715
$ws.Range("K136").Value = 8158.7145
$ws.Range("M136").Value = -5608.7145

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1610
$ws.Range("I64").Value = 1301
$ws.Range("J64").Value = 1816
$ws.Range("K64").Value = 1301
$ws.Range("L64").Value = 1816
$ws.Range("M64").Value = -1076
$ws.Range("N64").Value = -2266
$ws.Range("H67").Value = 1610
$ws.Range("I67").Value = 1301
$ws.Range("J67").Value = 1816
$ws.Range("K67").Value = 1301
$ws.Range("L67").Value = 1816
$ws.Range("M67").Value = -521
$ws.Range("N67").Value = -3376

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 624.3333
$ws.Range("I16").Value = 449.25
$ws.Range("J16").Value = 974.5
$ws.Range("K16").Value = 449.25
$ws.Range("L16").Value = 974.5
$ws.Range("M16").Value = -162.25
$ws.Range("N16").Value = -1548.5
$ws.Range("H31").Value = 4630.3
$ws.Range("I31").Value = 3787.35
$ws.Range("J31").Value = 5473.25
$ws.Range("K31").Value = 3787.35
$ws.Range("L31").Value = 5473.25
$ws.Range("M31").Value = -3492.35
$ws.Range("N31").Value = -6063.25
$ws.Range("H34").Value = 4630.3
$ws.Range("I34").Value = 3787.35
$ws.Range("J34").Value = 5473.25
$ws.Range("K34").Value = 3787.35
$ws.Range("L34").Value = 5473.25
$ws.Range("M34").Value = -3585.35
$ws.Range("N34").Value = -5877.25
$ws.Range("H113").Value = 624.3333
$ws.Range("I113").Value = 449.25
$ws.Range("J113").Value = 974.5
$ws.Range("K113").Value = 449.25
$ws.Range("L113").Value = 974.5
$ws.Range("M113").Value = 1720.75
$ws.Range("N113").Value = -5314.5
$ws.Range("H134").Value = 2893.1482
$ws.Range("I134").Value = 2256.6667
$ws.Range("K134").Value = 6770.000100000001
$ws.Range("M134").Value = -4235.000100000001

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 184.4
$ws.Range("I2").Value = 170.11111
$ws.Range("J2").Value = 313
$ws.Range("K2").Value = 1020.66666
$ws.Range("L2").Value = 1878
$ws.Range("M2").Value = -907.66666
$ws.Range("N2").Value = -2104
$ws.Range("H4").Value = 3000456.8
$ws.Range("I4").Value = 4200380
$ws.Range("J4").Value = 649.5
$ws.Range("K4").Value = 12601140
$ws.Range("L4").Value = 1948.5
$ws.Range("M4").Value = -12601028
$ws.Range("N4").Value = -2172.5
$ws.Range("H6").Value = 38.25
$ws.Range("I6").Value = 38.25
$ws.Range("K6").Value = 114.75
$ws.Range("M6").Value = -1.75
$ws.Range("H7").Value = 100000000
$ws.Range("I7").Value = 100000000
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 300000000
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -299999888
$ws.Range("N7").ClearContents()   # column no longer populated
$ws.Range("J31").Value = 500
$ws.Range("L31").Value = 1500
$ws.Range("N31").Value = -2076
$ws.Range("H50").Value = 165.28572
$ws.Range("I50").Value = 192
$ws.Range("K50").Value = 576
$ws.Range("M50").Value = -95
$ws.Range("H53").Value = 165.28572
$ws.Range("I53").Value = 192
$ws.Range("K53").Value = 576
$ws.Range("M53").Value = -95
$ws.Range("H114").Value = 937.6667
$ws.Range("I114").Value = 456.5
$ws.Range("J114").Value = 1900
$ws.Range("K114").Value = 1369.5
$ws.Range("L114").Value = 5700
$ws.Range("M114").Value = 1884.5
$ws.Range("N114").Value = -12208
$ws.Range("H140").Value = 4417.8667
$ws.Range("I140").Value = 3559.1538
$ws.Range("K140").Value = 10677.4614
$ws.Range("M140").Value = -5497.4614

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5672.706
$ws.Range("J80").Value = 6117.769
$ws.Range("L80").Value = 6117.769
$ws.Range("N80").Value = -8113.769
$ws.Range("H83").Value = 5672.706
$ws.Range("J83").Value = 6117.769
$ws.Range("L83").Value = 30588.845
$ws.Range("N83").Value = -40572.845
$ws.Range("H132").Value = 3553.8696
$ws.Range("I132").Value = 2485.2666
$ws.Range("K132").Value = 7455.7998
$ws.Range("M132").Value = -4925.7998

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3333.5557
$ws.Range("I81").Value = 2999.8572
$ws.Range("K81").Value = 5999.7144
$ws.Range("M81").Value = -4938.7144
$ws.Range("H84").Value = 3333.5557
$ws.Range("I84").Value = 2999.8572
$ws.Range("K84").Value = 29998.572
$ws.Range("M84").Value = -24694.572
$ws.Range("H132").Value = 2511.55
$ws.Range("I132").Value = 2015.6666
$ws.Range("K132").Value = 6046.9998
$ws.Range("M132").Value = -3516.9998
$ws.Range("H136").Value = 79618.08
$ws.Range("I136").Value = 1704.6
$ws.Range("K136").Value = 5113.799999999999
$ws.Range("M136").Value = -2563.799999999999

Write-Output "Applied scheduled Sheets refresh."